$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each target cell to Text format before writing, then clear the
# number-format override so the saved file has no stray style index —
# this avoids Excel auto-converting numeric-looking strings (e.g. "0.9999",
# "1.000") into real numbers while still matching the original (unstyled)
# plain-text cells.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '31.214.11'
Set-TextValue 'E2' '  +2.66%  '
Set-TextValue 'D3' '1.998.64'
Set-TextValue 'E3' '  +6.71%  '
Set-TextValue 'D4' '0.9999'
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '0.7816'
Set-TextValue 'E5' '  +65.59%  '
Set-TextValue 'D6' '257.04'
Set-TextValue 'E6' '  +4.68%  '
Set-TextValue 'D7' '0.9989'
Set-TextValue 'E7' '  -0.09%  '
Set-TextValue 'D8' '0.3530'
Set-TextValue 'E8' '  +23.12%  '
Set-TextValue 'D9' '29.29'
Set-TextValue 'E9' '  +34.13%  '
Set-TextValue 'B10' 'Dogecoin'
Set-TextValue 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D10' '0.07037'
Set-TextValue 'E10' '  +8.36%  '
Set-TextValue 'B11' 'Polygon'
Set-TextValue 'C11' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D11' '0.8672'
Set-TextValue 'E11' '  +18.97%  '
Set-TextValue 'B12' 'TRON'
Set-TextValue 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D12' '0.08205'
Set-TextValue 'E12' '  +5.09%  '
Set-TextValue 'B13' 'WrappedEther'
Set-TextValue 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D13' '1.998.13'
Set-TextValue 'E13' '  +6.73%  '
Set-TextValue 'B14' 'Litecoin'
Set-TextValue 'C14' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D14' '100.92'
Set-TextValue 'E14' '  +0.39%  '
Set-TextValue 'B15' 'Polkadot'
Set-TextValue 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '5.596'
Set-TextValue 'E15' '  +8.30%  '
Set-TextValue 'B16' 'Avalanche'
Set-TextValue 'C16' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D16' '15.59'
Set-TextValue 'E16' '  +19.12%  '
Set-TextValue 'B17' 'BitcoinCash'
Set-TextValue 'C17' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D17' '274.02'
Set-TextValue 'E17' '  -3.43%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '31.217.01'
Set-TextValue 'E18' '  +2.73%  '
Set-TextValue 'B19' 'Uniswap'
Set-TextValue 'C19' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D19' '5.957'
Set-TextValue 'E19' '  +11.83%  '
Set-TextValue 'B20' 'ShibaInu'
Set-TextValue 'C20' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D20' '0.000007944'
Set-TextValue 'E20' '  +6.13%  '
Set-TextValue 'B21' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C21' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D21' '2.263.73'
Set-TextValue 'E21' '  +7.10%  '
Set-TextValue 'B22' 'Dai'
Set-TextValue 'C22' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D22' '0.9987'
Set-TextValue 'E22' '  -0.13%  '
Set-TextValue 'B23' 'BinanceUSD'
Set-TextValue 'C23' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D23' '1.000'
Set-TextValue 'E23' '  +0.03%  '
Set-TextValue 'B24' 'Chainlink'
Set-TextValue 'C24' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D24' '7.140'
Set-TextValue 'E24' '  +12.78%  '
Set-TextValue 'B25' 'Cosmos'
Set-TextValue 'C25' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D25' '10.09'
Set-TextValue 'E25' '  +11.56%  '
Set-TextValue 'B26' 'Monero'
Set-TextValue 'C26' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D26' '164.57'
Set-TextValue 'E26' '  +1.57%  '
Set-TextValue 'D27' '0.1484'
Set-TextValue 'E27' '  +53.34%  '
Set-TextValue 'B28' 'EthereumClassic'
Set-TextValue 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D28' '19.99'
Set-TextValue 'E28' '  +5.43%  '
Set-TextValue 'B29' 'LidoDAOToken'
Set-TextValue 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D29' '2.378'
Set-TextValue 'E29' '  +25.45%  '
Set-TextValue 'B30' 'PancakeSwap'
Set-TextValue 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '1.609'
Set-TextValue 'E30' '  +7.86%  '
Set-TextValue 'B31' 'Filecoin'
Set-TextValue 'C31' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D31' '4.618'
Set-TextValue 'E31' '  +9.27%  '
Set-TextValue 'B32' 'Toncoin'
Set-TextValue 'C32' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D32' '1.362'
Set-TextValue 'E32' '  +3.08%  '
Set-TextValue 'B33' 'InternetComputer(DFINITY)'
Set-TextValue 'C33' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D33' '4.450'
Set-TextValue 'E33' '  +7.13%  '
Set-TextValue 'B34' 'Hedera'
Set-TextValue 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D34' '0.05220'
Set-TextValue 'E34' '  +8.40%  '
Set-TextValue 'B35' 'ARBITRUM'
Set-TextValue 'C35' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D35' '1.229'
Set-TextValue 'E35' '  +9.26%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.7788'
Set-TextValue 'E36' '  +12.72%  '
Set-TextValue 'B37' 'HuobiToken'
Set-TextValue 'C37' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D37' '2.801'
Set-TextValue 'E37' '  +2.14%  '
Set-TextValue 'B38' 'VeChain'
Set-TextValue 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D38' '0.02010'
Set-TextValue 'E38' '  +6.13%  '
Set-TextValue 'B39' 'MXToken'
Set-TextValue 'C39' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D39' '2.919'
Set-TextValue 'E39' '  +2.66%  '
Set-TextValue 'B40' 'FraxShare'
Set-TextValue 'C40' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D40' '6.745'
Set-TextValue 'E40' '  +7.02%  '
Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '79.74'
Set-TextValue 'E41' '  +4.88%  '
Set-TextValue 'B42' 'RenderToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D42' '2.162'
Set-TextValue 'E42' '  +10.75%  '
Set-TextValue 'D43' '0.4716'
Set-TextValue 'E43' '  +11.91%  '
Set-TextValue 'B44' 'Quant'
Set-TextValue 'C44' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D44' '106.52'
Set-TextValue 'E44' '  +5.54%  '
Set-TextValue 'B45' 'TrustWalletToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D45' '0.8505'
Set-TextValue 'E45' '  +3.03%  '
Set-TextValue 'B46' 'PaxDollar'
Set-TextValue 'C46' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D46' '0.9993'
Set-TextValue 'E46' '  +0.02%  '
Set-TextValue 'B47' 'Aptos'
Set-TextValue 'C47' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D47' '7.735'
Set-TextValue 'E47' '  +10.25%  '
Set-TextValue 'B48' 'EnergySwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '9.958'
Set-TextValue 'E48' '  +1.69%  '
Set-TextValue 'B49' 'Decentraland'
Set-TextValue 'C49' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D49' '0.4339'
Set-TextValue 'E49' '  +10.70%  '
Set-TextValue 'B50' 'Elrond'
Set-TextValue 'C50' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D50' '36.83'
Set-TextValue 'E50' '  +5.23%  '
Set-TextValue 'B51' 'NEARProtocol'
Set-TextValue 'C51' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D51' '1.513'
Set-TextValue 'E51' '  +13.81%  '
